# Apply the profit/price data updates scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1298.125
$ws.Range("I19").Value = 1710.6666
$ws.Range("K19").Value = 1710.6666
$ws.Range("M19").Value = -1535.6666

$ws.Range("H33").Value = 149.5
$ws.Range("I33").Value = 149.5
$ws.Range("K33").Value = 149.5
$ws.Range("M33").Value = 79.5

$ws.Range("H96").Value = 599.6
$ws.Range("I96").Value = 500
$ws.Range("J96").Value = 749
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2247
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -4993

$ws.Range("H127").Value = 4498
$ws.Range("J127").Value = 4498
$ws.Range("L127").Value = 13494
$ws.Range("N127").Value = -23414

$ws.Range("H131").Value = 3783.3333
$ws.Range("I131").Value = 675
$ws.Range("K131").Value = 2025
$ws.Range("M131").Value = 3015

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 545
$ws.Range("I141").Value = 545
$ws.Range("K141").Value = 1635
$ws.Range("M141").Value = 3545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1458.75
$ws.Range("I61").Value = 1458.75
$ws.Range("K61").Value = 1458.75
$ws.Range("M61").Value = -1246.75

$ws.Range("H132").Value = 1011.5
$ws.Range("I132").Value = 1011.5
$ws.Range("K132").Value = 3034.5
$ws.Range("M132").Value = -504.5

$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws.Range("H136").Value = 1458.75
$ws.Range("I136").Value = 1458.75
$ws.Range("K136").Value = 4376.25
$ws.Range("M136").Value = -1826.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2098.2
$ws.Range("J20").Value = 2500
$ws.Range("L20").Value = 2500
$ws.Range("N20").Value = -2994

$ws.Range("H110").Value = 49995
$ws.Range("J110").Value = 49995
$ws.Range("L110").Value = 49995
$ws.Range("N110").Value = -58175

$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("K10").Value = 200
$ws.Range("M10").Value = -61

$ws.Range("H31").Value = 3874
$ws.Range("I31").Value = 4999
$ws.Range("J31").Value = 3499
$ws.Range("K31").Value = 4999
$ws.Range("L31").Value = 3499
$ws.Range("M31").Value = -4704
$ws.Range("N31").Value = -4089

$ws.Range("H34").Value = 3874
$ws.Range("I34").Value = 4999
$ws.Range("J34").Value = 3499
$ws.Range("K34").Value = 4999
$ws.Range("L34").Value = 3499
$ws.Range("M34").Value = -4797
$ws.Range("N34").Value = -3903

$ws.Range("H58").Value = 3650.75
$ws.Range("I58").Value = 3001
$ws.Range("K58").Value = 3001
$ws.Range("M58").Value = -2798

$ws.Range("H132").Value = 1879.8
$ws.Range("I132").Value = 1533.1111
$ws.Range("K132").Value = 4599.3333
$ws.Range("M132").Value = -2069.3333

$ws.Range("H134").Value = 2937.2942
$ws.Range("I134").Value = 2363.6155
$ws.Range("J134").Value = 4801.75
$ws.Range("K134").Value = 7090.8465
$ws.Range("L134").Value = 14405.25
$ws.Range("M134").Value = -4555.8465
$ws.Range("N134").Value = -19475.25

$ws.Range("H136").Value = 3650.75
$ws.Range("I136").Value = 3001
$ws.Range("K136").Value = 9003
$ws.Range("M136").Value = -6453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 59.25
$ws.Range("I6").Value = 67.333336
$ws.Range("K6").Value = 202.000008
$ws.Range("M6").Value = -89.00000800000001

$ws.Range("H15").Value = 1099.5
$ws.Range("J15").Value = 1099.5
$ws.Range("L15").Value = 3298.5
$ws.Range("N15").Value = -3578.5

$ws.Range("H104").Value = 5899
$ws.Range("J104").Value = 5899
$ws.Range("L104").Value = 17697
$ws.Range("N104").Value = -22939

$ws.Range("H131").Value = 2000
$ws.Range("J131").Value = 2500
$ws.Range("L131").Value = 7500
$ws.Range("N131").Value = -17580

$ws.Range("H134").Value = 11333
$ws.Range("I134").Value = 11333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 33999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -28929
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4837.4165
$ws.Range("I102").Value = 4504.9
$ws.Range("J102").Value = 6500
$ws.Range("K102").Value = 4504.9
$ws.Range("L102").Value = 6500
$ws.Range("M102").Value = -2882.9
$ws.Range("N102").Value = -9744

$ws.Range("H113").Value = 1249.5
$ws.Range("I113").Value = 1249.5
$ws.Range("K113").Value = 1249.5
$ws.Range("M113").Value = 920.5

$ws.Range("H126").Value = 10505.5
$ws.Range("I126").Value = 10505.5
$ws.Range("K126").Value = 31516.5
$ws.Range("M126").Value = -29046.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214

$ws.Range("H55").Value = 2137.7144
$ws.Range("J55").Value = 2593.4
$ws.Range("L55").Value = 2593.4
$ws.Range("N55").Value = -2939.4

$ws.Range("H101").Value = 15739
$ws.Range("J101").Value = 15739
$ws.Range("L101").Value = 15739
$ws.Range("N101").Value = -22229

$ws.Range("H132").Value = 6408.3335
$ws.Range("I132").Value = 6612.5
$ws.Range("K132").Value = 19837.5
$ws.Range("M132").Value = -17307.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 3695.7778
$ws.Range("I132").Value = 3220.25
$ws.Range("K132").Value = 9660.75
$ws.Range("M132").Value = -7130.75
